# "Generate Report for Handoff" — mark b.md as ready for handoff by
# updating the Overview, zh-cn and de-de sheets with the new handoff
# file names / timestamps, and flipping status from
# "Handed back: in sync with en-US" to "Ready for handoff".

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) status columns for zh-cn / de-de
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $readyForHandoff
$overview.Range("C3").Value = $readyForHandoff

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("B3").Value = $readyForHandoff
$zhcn.Range("C3").Value = $zhcnHandoffFile
$zhcn.Range("D3").Value = "2016-03-07 02:15:09"

foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = $zhcnHandoffFile
    }
}

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dedeHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("B3").Value = $readyForHandoff
$dede.Range("C3").Value = $dedeHandoffFile
$dede.Range("D3").Value = "2016-03-07 02:15:19"

foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = $dedeHandoffFile
    }
}
